$d = $word.ActiveDocument

$guideText = "This document serves as a guide for a project designed to demonstrate the process of generating PDF documents from Markdown using Quarto and LaTeX. Additionally, it incorporates the use of GitHub Actions to automate the generation process and GitHub Releases for storing the final output. This repository has been structured as a GitHub Template, allowing it to be easily used as a starting point for any book project or documentation endeavor."

# Paragraph 7: "2. Overview of Technologies" (Heading1) -> "2. Commands"
$d.Paragraphs.Item(7).Range.Text = "2. Commands"

# Paragraph 8: intro blurb (FirstParagraph) -> repeat the guide text
$d.Paragraphs.Item(8).Range.Text = $guideText

# Paragraph 9: "2.0.1 Quarto" (Heading3) -> "3. Containers" (Heading1)
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "3. Containers"
$p9.Style = "Heading 1"

# Paragraph 10: Quarto description (FirstParagraph) -> repeat the guide text again
$d.Paragraphs.Item(10).Range.Text = $guideText

# Remove everything from paragraph 11 through the end of the document
# (LaTeX/GitHub Actions/GitHub Release subsections, Project Structure,
#  Workflow, GitHub Template, and Conclusion sections).
$startRange = $d.Paragraphs.Item(11).Range.Start
$endRange = $d.Content.End
$d.Range($startRange, $endRange).Delete()
